# Generate Report for Handback
# The localized file "cc0121d2-f0d0-4fec-8829-8fd6ee143b4b.md" has been handed
# back and is now in sync with en-US, so update its status on every sheet and
# refresh the "Latest Handback DateTime" on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 is the cc0121d2... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: row 3 is the cc0121d2... file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $status
$wsZhCn.Range("G3").Value = "2016-03-10 07:10:43"

# --- de-de sheet: row 3 is the cc0121d2... file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $status
$wsDeDe.Range("G3").Value = "2016-03-10 07:10:55"
